$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set CommentRow (column K) values to 1 for rows 2 through 27
$ws.Range("K2:K27").Value = 1
